$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "ROMERO CHANAME YOSSELY TRINIDAD",
    "ZAVALETA MANAY JORGE LUIS",
    "HIDALGO CUBAS LUISA YVONE",
    "SENADOR ARBOLEDA GIANCARLOS EXEBIO",
    "ZEVALLOS PACHECO ZOILA XIMENA",
    "CONTRERAS VALDERRAMA JULIA ALEJANDRA",
    "FERNANDEZ VALDERAS ERNESTO ALI",
    "VALLE MAGALLAN EDUAR",
    "CAMACHO LINARES JUDITH ARLETT",
    "HUMPIRE CASTILLO IRWIN DEIMER",
    "SEVERINO AVALOS MARJORIE ISABEL",
    "BALLENA ESQUÉN ASTRID CAROLINA",
    "GUTIERREZ CARLOS TERESA DE JESUS"
)

$values = @(42, 40, 39, 39, 35, 34, 33, 32, 30, 30, 29, 21, 16)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
